$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.379393333333333
$ws.Range("H2").Value = 4.13818
$ws.Range("I2").Value = 0.2804878676989906
$ws.Range("J2").Value = 0.2804878676989906
$ws.Range("M2").Value = 0.5550926666666666
$ws.Range("N2").Value = 1.665278
$ws.Range("O2").Value = 0.1208967663154349
$ws.Range("P2").Value = 0.1208967663154349
$ws.Range("Q2").Value = 0.7656911237822221
$ws.Range("R2").Value = 6.891220114039999
$ws.Range("S2").Value = 0.03391007619551949
$ws.Range("T2").Value = 0.03391007619551948
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.379393333333333
$ws.Range("H3").Value = 4.13818
$ws.Range("I3").Value = 0.2804878676989906
$ws.Range("J3").Value = 0.2804878676989906
$ws.Range("O3").Value = 0.7377399926530269
$ws.Range("P3").Value = 0.7377399926530268
$ws.Range("Q3").Value = 4.672424095775556
$ws.Range("R3").Value = 42.05181686198
$ws.Range("S3").Value = 0.2069271174555165
$ws.Range("T3").Value = 0.2069271174555165
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.379393333333333
$ws.Range("H4").Value = 4.13818
$ws.Range("I4").Value = 0.2804878676989906
$ws.Range("J4").Value = 0.2804878676989906
$ws.Range("M4").Value = 0.5311786666666667
$ws.Range("N4").Value = 1.593536
$ws.Range("O4").Value = 0.1156884012202364
$ws.Range("P4").Value = 0.1156884012202364
$ws.Range("Q4").Value = 0.7327043116088889
$ws.Range("R4").Value = 6.59433880448
$ws.Range("S4").Value = 0.03244919297576942
$ws.Range("T4").Value = 0.03244919297576942
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.379393333333333
$ws.Range("H5").Value = 4.13818
$ws.Range("I5").Value = 0.2804878676989906
$ws.Range("J5").Value = 0.2804878676989906
$ws.Range("M5").Value = 0.117885
$ws.Range("N5").Value = 0.353655
$ws.Range("O5").Value = 0.02567483981130185
$ws.Range("P5").Value = 0.02567483981130185
$ws.Range("Q5").Value = 0.1626097831
$ws.Range("R5").Value = 1.4634880479
$ws.Range("S5").Value = 0.007201481072185212
$ws.Range("T5").Value = 0.007201481072185211
$ws.Range("I6").Value = 0.4808449624319857
$ws.Range("J6").Value = 0.4808449624319857
$ws.Range("M6").Value = 0.5550926666666666
$ws.Range("N6").Value = 1.665278
$ws.Range("O6").Value = 0.1208967663154349
$ws.Range("P6").Value = 0.1208967663154349
$ws.Range("Q6").Value = 1.312636880411111
$ws.Range("R6").Value = 11.8137319237
$ws.Range("S6").Value = 0.05813260105709383
$ws.Range("T6").Value = 0.05813260105709382
$ws.Range("I7").Value = 0.4808449624319857
$ws.Range("J7").Value = 0.4808449624319857
$ws.Range("O7").Value = 0.7377399926530269
$ws.Range("P7").Value = 0.7377399926530268
$ws.Range("S7").Value = 0.3547385590518181
$ws.Range("T7").Value = 0.3547385590518181
$ws.Range("I8").Value = 0.4808449624319857
$ws.Range("J8").Value = 0.4808449624319857
$ws.Range("M8").Value = 0.5311786666666667
$ws.Range("N8").Value = 1.593536
$ws.Range("O8").Value = 0.1156884012202364
$ws.Range("P8").Value = 0.1156884012202364
$ws.Range("Q8").Value = 1.256087046044445
$ws.Range("R8").Value = 11.3047834144
$ws.Range("S8").Value = 0.05562818493856107
$ws.Range("T8").Value = 0.05562818493856106
$ws.Range("I9").Value = 0.4808449624319857
$ws.Range("J9").Value = 0.4808449624319857
$ws.Range("M9").Value = 0.117885
$ws.Range("N9").Value = 0.353655
$ws.Range("O9").Value = 0.02567483981130185
$ws.Range("P9").Value = 0.02567483981130185
$ws.Range("Q9").Value = 0.27876462425
$ws.Range("R9").Value = 2.50888161825
$ws.Range("S9").Value = 0.01234561738451269
$ws.Range("T9").Value = 0.01234561738451269
$ws.Range("G10").Value = 1.173726
$ws.Range("H10").Value = 3.521178
$ws.Range("I10").Value = 0.2386671698690237
$ws.Range("J10").Value = 0.2386671698690237
$ws.Range("M10").Value = 0.5550926666666666
$ws.Range("N10").Value = 1.665278
$ws.Range("O10").Value = 0.1208967663154349
$ws.Range("P10").Value = 0.1208967663154349
$ws.Range("Q10").Value = 0.6515266952759999
$ws.Range("R10").Value = 5.863740257483999
$ws.Range("S10").Value = 0.02885408906282156
$ws.Range("T10").Value = 0.02885408906282155
$ws.Range("G11").Value = 1.173726
$ws.Range("H11").Value = 3.521178
$ws.Range("I11").Value = 0.2386671698690237
$ws.Range("J11").Value = 0.2386671698690237
$ws.Range("O11").Value = 0.7377399926530269
$ws.Range("P11").Value = 0.7377399926530268
$ws.Range("Q11").Value = 3.975766383462
$ws.Range("R11").Value = 35.781897451158
$ws.Range("S11").Value = 0.1760743161456923
$ws.Range("T11").Value = 0.1760743161456922
$ws.Range("G12").Value = 1.173726
$ws.Range("H12").Value = 3.521178
$ws.Range("I12").Value = 0.2386671698690237
$ws.Range("J12").Value = 0.2386671698690237
$ws.Range("M12").Value = 0.5311786666666667
$ws.Range("N12").Value = 1.593536
$ws.Range("O12").Value = 0.1156884012202364
$ws.Range("P12").Value = 0.1156884012202364
$ws.Range("Q12").Value = 0.623458211712
$ws.Range("R12").Value = 5.611123905408
$ws.Range("S12").Value = 0.02761102330590593
$ws.Range("T12").Value = 0.02761102330590593
$ws.Range("G13").Value = 1.173726
$ws.Range("H13").Value = 3.521178
$ws.Range("I13").Value = 0.2386671698690237
$ws.Range("J13").Value = 0.2386671698690237
$ws.Range("M13").Value = 0.117885
$ws.Range("N13").Value = 0.353655
$ws.Range("O13").Value = 0.02567483981130185
$ws.Range("P13").Value = 0.02567483981130185
$ws.Range("Q13").Value = 0.13836468951
$ws.Range("R13").Value = 1.24528220559
$ws.Range("S13").Value = 0.006127741354603952
$ws.Range("T13").Value = 0.006127741354603951
